# Updates the cryptos price table (columns D = Price, E = Volume(1h))
# for rows 2-51, per the scraped-data refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.785.36"
$ws.Range("E2").Value = "  +2.50%  "
$ws.Range("D3").Value = "'2.422.48"
$ws.Range("E3").Value = "  +3.24%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'553.50"
$ws.Range("E5").Value = "  +1.93%  "
$ws.Range("D6").Value = "'137.02"
$ws.Range("E6").Value = "  +1.32%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +1.15%  "
$ws.Range("E9").Value = "  +5.16%  "
$ws.Range("E10").Value = "  +2.52%  "
$ws.Range("E11").Value = "  +1.77%  "
$ws.Range("E12").Value = "  -2.25%  "
$ws.Range("E13").Value = "  +3.43%  "
$ws.Range("D14").Value = "'2.849.67"
$ws.Range("E14").Value = "  +3.13%  "
$ws.Range("D15").Value = "'59.670.99"
$ws.Range("E15").Value = "  +2.43%  "
$ws.Range("D16").Value = "'0.0000139"
$ws.Range("E16").Value = "  +4.20%  "
$ws.Range("D17").Value = "'2.425.35"
$ws.Range("E17").Value = "  +3.66%  "
$ws.Range("D18").Value = "'11.31"
$ws.Range("E18").Value = "  +5.50%  "
$ws.Range("E19").Value = "  +4.47%  "
$ws.Range("D20").Value = "'336.43"
$ws.Range("E20").Value = "  +0.81%  "
$ws.Range("D21").Value = "'6.97"
$ws.Range("E21").Value = "  +4.64%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "'64.62"
$ws.Range("E23").Value = "  +2.95%  "
$ws.Range("E24").Value = "  +0.94%  "
$ws.Range("E25").Value = "  +0.32%  "
$ws.Range("D26").Value = "'0.997"
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("D27").Value = "'1.38"
$ws.Range("E27").Value = "  -2.17%  "
$ws.Range("D28").Value = "'0.0₃0783"
$ws.Range("E28").Value = "  +6.13%  "
$ws.Range("D29").Value = "'1.80"
$ws.Range("E29").Value = "  +2.47%  "
$ws.Range("D30").Value = "'170.66"
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("E31").Value = "  +2.59%  "
$ws.Range("D32").Value = "'18.73"
$ws.Range("E32").Value = "  +1.65%  "
$ws.Range("D33").Value = "'1.02"
$ws.Range("E33").Value = "  -0.43%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").Value = "'1.31"
$ws.Range("E35").Value = "  +5.25%  "
$ws.Range("D36").Value = "'4.30"
$ws.Range("E36").Value = "  +0.66%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.40%  "
$ws.Range("D38").Value = "'1.64"
$ws.Range("E38").Value = "  -0.52%  "
$ws.Range("D39").Value = "'40.15"
$ws.Range("E39").Value = "  +2.57%  "
$ws.Range("D40").Value = "'0.418"
$ws.Range("E40").Value = "  +11.08%  "
$ws.Range("D41").Value = "'305.83"
$ws.Range("E41").Value = "  +5.93%  "
$ws.Range("E42").Value = "  +2.95%  "
$ws.Range("D43").Value = "'142.51"
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").Value = "'0.0963"
$ws.Range("E44").Value = "  +2.72%  "
$ws.Range("D45").Value = "'0.0525"
$ws.Range("E45").Value = "  +4.23%  "
$ws.Range("E46").Value = "  +1.67%  "
$ws.Range("D47").Value = "'19.09"
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("D48").Value = "'0.404"
$ws.Range("E48").Value = "  +4.74%  "
$ws.Range("D49").Value = "'0.0225"
$ws.Range("E49").Value = "  +2.91%  "
$ws.Range("E50").Value = "  -0.31%  "
$ws.Range("D51").Value = "'1.61"
$ws.Range("E51").Value = "  +4.67%  "
